$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 82.987681
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 1.150126
$ws.Range("N2").Value = 3.450378
$ws.Range("O2").Value = 0.3945373201707233
$ws.Range("P2").Value = 0.3945373201707233
$ws.Range("Q2").Value = 95.44628959780599
$ws.Range("R2").Value = 859.016606380254
$ws.Range("S2").Value = 0.1771276922596378
$ws.Range("T2").Value = 0.1771276922596378

$ws.Range("G3").Value = 82.987681
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.2405487790236168
$ws.Range("P3").Value = 0.2405487790236168
$ws.Range("Q3").Value = 58.19345154762999
$ws.Range("R3").Value = 523.74106392867
$ws.Range("S3").Value = 0.1079944733387697
$ws.Range("T3").Value = 0.1079944733387697

$ws.Range("G4").Value = 82.987681
$ws.Range("H4").Value = 248.963043
$ws.Range("I4").Value = 0.4489504115427952
$ws.Range("J4").Value = 0.4489504115427952
$ws.Range("O4").Value = 0.3649139008056599
$ws.Range("P4").Value = 0.3649139008056599
$ws.Range("Q4").Value = 88.27980541737
$ws.Range("R4").Value = 794.51824875633
$ws.Range("S4").Value = 0.1638282459443878
$ws.Range("T4").Value = 0.1638282459443878

$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("M5").Value = 1.150126
$ws.Range("N5").Value = 3.450378
$ws.Range("O5").Value = 0.3945373201707233
$ws.Range("P5").Value = 0.3945373201707233
$ws.Range("Q5").Value = 72.61963344758934
$ws.Range("R5").Value = 653.5767010283039
$ws.Range("S5").Value = 0.1347663501589692
$ws.Range("T5").Value = 0.1347663501589692

$ws.Range("G6").Value = 63.14058933333333
$ws.Range("I6").Value = 0.3415807409566563
$ws.Range("J6").Value = 0.3415807409566563
$ws.Range("O6").Value = 0.2405487790236168
$ws.Range("P6").Value = 0.2405487790236168
$ws.Range("R6").Value = 398.4846791239199
$ws.Range("S6").Value = 0.08216683017510601
$ws.Range("T6").Value = 0.08216683017510602

$ws.Range("G7").Value = 63.14058933333333
$ws.Range("I7").Value = 0.3415807409566563
$ws.Range("J7").Value = 0.3415807409566563
$ws.Range("O7").Value = 0.3649139008056599
$ws.Range("P7").Value = 0.3649139008056599
$ws.Range("S7").Value = 0.1246475606225811
$ws.Range("T7").Value = 0.1246475606225811

$ws.Range("I8").Value = 0.2094688475005485
$ws.Range("J8").Value = 0.2094688475005485
$ws.Range("M8").Value = 1.150126
$ws.Range("N8").Value = 3.450378
$ws.Range("O8").Value = 0.3945373201707233
$ws.Range("P8").Value = 0.3945373201707233
$ws.Range("Q8").Value = 44.53281201269201
$ws.Range("R8").Value = 400.795308114228
$ws.Range("S8").Value = 0.08264327775211631
$ws.Range("T8").Value = 0.08264327775211631

$ws.Range("I9").Value = 0.2094688475005485
$ws.Range("J9").Value = 0.2094688475005485
$ws.Range("O9").Value = 0.2405487790236168
$ws.Range("P9").Value = 0.2405487790236168
$ws.Range("Q9").Value = 27.15158492866
$ws.Range("S9").Value = 0.05038747550974112
$ws.Range("T9").Value = 0.05038747550974113

$ws.Range("I10").Value = 0.2094688475005485
$ws.Range("J10").Value = 0.2094688475005485
$ws.Range("O10").Value = 0.3649139008056599
$ws.Range("P10").Value = 0.3649139008056599
$ws.Range("Q10").Value = 41.18911270134001
$ws.Range("S10").Value = 0.07643809423869105
$ws.Range("T10").Value = 0.07643809423869105
